$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the status of the existing "BasicDataStructure" row (row 3) ---
# "On progress" -> "Done"
$ws.Range("B3").Value = "Done"

# --- Add a new tracked exercise row (row 4): FoodOrderingSystem / Ongoing ---
$ws.Range("A4").Value = "FoodOrderingSystem"
$ws.Range("B4").Value = "Ongoing"

# --- Grow the Excel table ("Table2") so it covers the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B4"))

# --- Move the active selection like the author left it (B5) ---
[void]$ws.Range("B5").Select()

# --- Conditional formatting ---
# 1) A 3-color scale (red/yellow/green) that was already on B2 -- re-added here
#    with Excel's default colors, matching the lower-priority (priority=2) rule.
$rng1 = $ws.Range("B2")
$cs1 = $rng1.FormatConditions.AddColorScale(3)

# 2) A 2-color scale (white -> green) across the whole data column B2:B4,
#    which takes the higher priority (priority=1) since it was added last.
$rng2 = $ws.Range("B2:B4")
$cs2 = $rng2.FormatConditions.AddColorScale(2)
$cs2.ColorScaleCriteria(1).FormatColor.Color = 16776444
$cs2.ColorScaleCriteria(2).FormatColor.Color = 8109667

# Fix up priorities/order to match: B2 rule = priority 2, B2:B4 rule = priority 1
$cs1.Priority = 2
$cs2.Priority = 1
